# Add the Drone ("드론 kit") purchase entry as a new row (row 8) on the
# 장부 (ledger) worksheet.
#
# Columns: B=날짜(date) C=사용자(user) D=항목(item) E=금액(amount)
#          F=영수증 제출(receipt submitted) G=입금 완료(payment complete) H=비고(note)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B8: purchase date 2017-06-07 (Excel serial 42893), formatted like a date
$ws.Range("B8").Value = 42893
$ws.Range("B8").NumberFormat = "mm-dd-yy"

# C8: purchaser / user
$ws.Range("C8").Value = "지무근"

# D8: item purchased
$ws.Range("D8").Value = "드론 kit"

# E8: amount spent
$ws.Range("E8").Value = 160430

# G8: payment complete marker
$ws.Range("G8").Value = "O"
